$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 3), mirroring the structure of row 2
$ws.Range("A3").Value = 42600.881053240744
$ws.Range("B3").Value = "Named"
$ws.Range("C3").Value = 12681
$ws.Range("D3").Value = 7883
$ws.Range("E3").Value = 452
$ws.Range("F3").Value = 67
$ws.Range("G3").Value = 45
$ws.Range("H3").Value = 59
$ws.Range("I3").Value = 40
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0

# Column A widens slightly to accommodate the new value
$ws.Columns.Item(1).ColumnWidth = 14
